$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 306
# from serial date 45181 to 45182 (2023-09-12 -> 2023-09-13).
$ws.Range("C2:C306").Value = 45182
